$d = $word.ActiveDocument

# 1. Insert the "Vi förväntar oss..." paragraph right after the "Nedan
#    presenteras..." paragraph near the top of the document.
$afterPara = $d.Paragraphs.Item(3)
$afterPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(4)
$newPara.Range.Text = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."

# 2. Remove the two trailing empty paragraphs and the duplicate "Vi
#    förväntar oss..." paragraph that used to sit at the very end of the
#    document (now redundant since the text was moved up in step 1).
$count = $d.Paragraphs.Count
$startPara = $d.Paragraphs.Item($count - 2)
$endPara = $d.Paragraphs.Item($count)
$rangeToDelete = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rangeToDelete.Delete()

# 3. Update the date shown in the first-page header from 2023-11-13 to
#    2023-11-14.
$d.Sections.Item(1).Headers.Item(2).Range.Find.Execute("2023-11-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-14", 2)
